# "updated the cap rate" -- switch the perpetuity-formula cap rate from
# (0.032-0.015) to (0.06-0.02) on the Total sheet, wire the Montreal/Toronto
# amortization sheets' B columns to the Total sheet via formula instead of
# pasted constants, and add the new "time to save a down payment" section.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("Total")
$montreal = $wb.Worksheets.Item("Montreal")
$toronto = $wb.Worksheets.Item("Toronto")

# --- Total sheet: updated cap rate -------------------------------------
$total.Range("P1").Formula = "=698*12/(0.06-0.02)"
$total.Range("P3").Formula = "=1194*12/(0.06-0.02)"

# --- Montreal sheet: B2:B26 now reference Total!$I$6 instead of a pasted value
for ($r = 2; $r -le 26; $r++) {
    $montreal.Range("B$r").Formula = "=Total!`$I`$6"
}

# --- Toronto sheet: B2:B26 now reference Total!$I$7 instead of a pasted value,
#     and pick up the currency number format (style index 1) used elsewhere
#     on this workbook
for ($r = 2; $r -le 26; $r++) {
    $cell = $toronto.Range("B$r")
    $cell.Formula = "=Total!`$I`$7"
    $cell.NumberFormat = "`"$`"#,##0.00;[Red]\-`"$`"#,##0.00"
}

# --- Total sheet: new "time needed to accumulate down payment" section ---
$total.Range("A19").Value = "Time needed to accumulate down payment"

$total.Range("A21").Value = "Toronto"
$total.Range("B21").Value = "assuming a salary of 100K then 60K disposable income. 5000 dollars a month. 1000 dollars expenses. 1194 rent 500 dollars incidents so => "
$total.Range("O21").Formula = "=5000-(1000+1194+500)"
$total.Range("Q21").Value = "time to save"
$total.Range("S21").Formula = "=L7/O21/12"
$total.Range("T21").Value = "years"

$total.Range("A22").Value = "Montreal"
$total.Range("B22").Value = "assuming a salary of 100K then 60K disposable income. 5000 dollars a month. 1000 dollars expenses. 698 rent 500 dollars incidents so => "
$total.Range("O22").Formula = "=5000-(1000+698+500)"
$total.Range("Q22").Value = "time to save"
$total.Range("S22").Formula = "=L6/O22/12"
$total.Range("T22").Value = "years"
$total.Range("V22").Value = "takes twice as long for toronto than montreal"

# --- view bookkeeping to match where the author left the selection -------
# (select on the non-active sheets first, finish on Total so it stays the
# sheet that's tabSelected when the file is saved, like the original)
$montreal.Range("F6").Select() | Out-Null
$toronto.Range("B2:B26").Select() | Out-Null
$total.Range("V23").Select() | Out-Null
